$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Regenerate report: update the "Latest Handoff Date/Datetime" for the
# 1d967de0-8e3d-42cb-8063-a153b49e1ad3 row (row 5) following a new handoff.
$wsOverview.Range("D5").Value = "2016-03-22 16:44:04"
$wsZhCn.Range("E5").Value = "2016-03-22 16:43:59"
$wsDeDe.Range("E5").Value = "2016-03-22 16:44:04"
